$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2553.182
$ws.Range("I9").Value = 271.25
$ws.Range("K9").Value = 271.25
$ws.Range("M9").Value = -102.25

$ws.Range("H33").Value = 355.33334
$ws.Range("I33").Value = 349.5
$ws.Range("J33").Value = 402
$ws.Range("K33").Value = 349.5
$ws.Range("L33").Value = 402
$ws.Range("M33").Value = -120.5
$ws.Range("N33").Value = -860

$ws.Range("H69").Value = 5333.3335

$ws.Range("H72").Value = 5333.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 95.75
$ws.Range("I97").Value = 95.75
$ws.Range("K97").Value = 95.75
$ws.Range("M97").Value = 400.25

$ws.Range("H102").Value = 3070.8572
$ws.Range("I102").Value = 1499.4
$ws.Range("K102").Value = 1499.4
$ws.Range("M102").Value = 122.5999999999999

$ws.Range("H132").Value = 3898
$ws.Range("I132").Value = 3898
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11694
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9164
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 62.333332
$ws.Range("I5").Value = 62.333332
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 62.333332
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 50.666668
$ws.Range("N5").ClearContents()

$ws.Range("H7").Value = 3167086.5
$ws.Range("I7").Value = 6333689.5
$ws.Range("J7").Value = 483.33334
$ws.Range("K7").Value = 6333689.5
$ws.Range("L7").Value = 483.33334
$ws.Range("M7").Value = -6333576.5
$ws.Range("N7").Value = -709.33334

$ws.Range("H22").Value = 558.7143
$ws.Range("I22").Value = 318.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 318.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -145.5
$ws.Range("N22").Value = -2346

$ws.Range("H105").Value = 3279.818
$ws.Range("I105").Value = 2634.875
$ws.Range("K105").Value = 2634.875
$ws.Range("M105").Value = -887.875

$ws.Range("H134").Value = 4697.4287
$ws.Range("I134").Value = 4577.6
$ws.Range("K134").Value = 13732.8
$ws.Range("M134").Value = -11197.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1457.8182
$ws.Range("I22").Value = 224
$ws.Range("J22").Value = 1732
$ws.Range("K22").Value = 224
$ws.Range("L22").Value = 1732
$ws.Range("M22").Value = 126
$ws.Range("N22").Value = -2432

$ws.Range("H58").Value = 7356.1
$ws.Range("I58").Value = 6508.143
$ws.Range("J58").Value = 9334.666999999999
$ws.Range("K58").Value = 6508.143
$ws.Range("L58").Value = 9334.666999999999
$ws.Range("M58").Value = -6305.143
$ws.Range("N58").Value = -9740.666999999999

$ws.Range("H132").Value = 9498.286
$ws.Range("I132").Value = 8318.799999999999
$ws.Range("K132").Value = 24956.4
$ws.Range("M132").Value = -22426.4

$ws.Range("H134").Value = 2793.4285
$ws.Range("I134").Value = 2384
$ws.Range("K134").Value = 7152
$ws.Range("M134").Value = -4617

$ws.Range("H136").Value = 7356.1
$ws.Range("I136").Value = 6508.143
$ws.Range("J136").Value = 9334.666999999999
$ws.Range("K136").Value = 19524.429
$ws.Range("L136").Value = 28004.001
$ws.Range("M136").Value = -16974.429
$ws.Range("N136").Value = -33104.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1081.4286
$ws.Range("I86").Value = 655.625
$ws.Range("J86").Value = 1649.1666
$ws.Range("K86").Value = 1966.875
$ws.Range("L86").Value = 4947.4998
$ws.Range("M86").Value = -780.875
$ws.Range("N86").Value = -7319.4998

$ws.Range("H89").Value = 1081.4286
$ws.Range("I89").Value = 655.625
$ws.Range("J89").Value = 1649.1666
$ws.Range("K89").Value = 5900.625
$ws.Range("L89").Value = 14842.4994
$ws.Range("M89").Value = 27.375
$ws.Range("N89").Value = -26698.4994

$ws.Range("H113").Value = 573.5
$ws.Range("J113").Value = 573.5
$ws.Range("L113").Value = 1720.5
$ws.Range("N113").Value = -6060.5

$ws.Range("H132").Value = 2062.2144
$ws.Range("I132").Value = 1135.4
$ws.Range("J132").Value = 2577.111
$ws.Range("K132").Value = 10218.6
$ws.Range("L132").Value = 23193.999
$ws.Range("M132").Value = -7688.6
$ws.Range("N132").Value = -28253.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1899.5
$ws.Range("I80").Value = 1899.5
$ws.Range("K80").Value = 1899.5
$ws.Range("M80").Value = -901.5

$ws.Range("H83").Value = 1899.5
$ws.Range("I83").Value = 1899.5
$ws.Range("K83").Value = 9497.5
$ws.Range("M83").Value = -4505.5

$ws.Range("H122").Value = 1780
$ws.Range("I122").Value = 1766
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 5298
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -2848
$ws.Range("N122").Value = -10450

$ws.Range("H132").Value = 3175.739
$ws.Range("I132").Value = 3138.2727
$ws.Range("K132").Value = 9414.8181
$ws.Range("M132").Value = -6884.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2499.5
$ws.Range("J22").Value = 2499.5
$ws.Range("L22").Value = 2499.5
$ws.Range("N22").Value = -3089.5

$ws.Range("H27").Value = 2499.5
$ws.Range("J27").Value = 2499.5
$ws.Range("L27").Value = 2499.5
$ws.Range("N27").Value = -2713.5

$ws.Range("H46").Value = 2376.4707
$ws.Range("J46").Value = 2900
$ws.Range("L46").Value = 2900
$ws.Range("N46").Value = -3276

$ws.Range("H68").Value = 9500
$ws.Range("J68").Value = 9500
$ws.Range("L68").Value = 9500
$ws.Range("N68").Value = -10998

$ws.Range("H71").Value = 9500
$ws.Range("J71").Value = 9500
$ws.Range("L71").Value = 47500
$ws.Range("N71").Value = -54988

$ws.Range("H101").Value = 33989.6
$ws.Range("J101").Value = 33989.6
$ws.Range("L101").Value = 33989.6
$ws.Range("N101").Value = -40479.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 25000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 25000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H136").Value = 1925.2858
$ws.Range("I136").Value = 1302.3572
$ws.Range("K136").Value = 3907.0716
$ws.Range("M136").Value = -1357.0716

Write-Host "All edits applied"
